# Major code refactoring #422
#
# Adds three new worksheets (SequenceWithAndSplit, NestedAndSplit,
# AndSplitWithLoop) to the workbook, each describing a small activity
# "layout" table (class / activityReference / name columns) similar in
# shape to the existing "StartWithAndSplit" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a single data row (columns A,B,C) into a worksheet,
# re-using the number formatting/style of an already formatted row by
# copying it first (this keeps the existing cellXfs / shared style ids
# intact instead of Excel minting brand new style entries).
# ---------------------------------------------------------------------
# NOTE: this PowerShell-ish engine does not bind named parameters
# (`-Foo bar`) reliably, so every helper below takes plain positional
# arguments only.
function Set-Row3 {
    param($Sheet, $RowNum, $A, $B, $C)

    $Sheet.Range("A$RowNum").Value = $A
    $Sheet.Range("B$RowNum").Value = $B
    $Sheet.Range("C$RowNum").Value = $C
}

function New-LayoutSheet {
    param($SheetName, $TotalRows)

    $template = $wb.Worksheets.Item("StartWithAndSplit")
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $after)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $SheetName

    $existingRows = 10

    # Grow the sheet (copy formats of the last template row downward) if
    # the new table needs more rows than the template already has.
    if ($TotalRows -gt $existingRows) {
        $newSheet.Range("A$existingRows`:C$existingRows").Copy()
        $newSheet.Range("A$($existingRows+1):C$TotalRows").PasteSpecial(-4122)
    }

    return $newSheet
}

# ---------------------------------------------------------------------
# Sheet: SequenceWithAndSplit
# ---------------------------------------------------------------------
$ws = New-LayoutSheet "SequenceWithAndSplit" 14

Set-Row3 $ws 1 "layout" $null $null
Set-Row3 $ws 2 "class" "activityReference" "name"
Set-Row3 $ws 3 "Elementary" "TestItem_First:0" "First"
Set-Row3 $ws 4 "AndSplit" "" ""
Set-Row3 $ws 5 "Block" "" ""
Set-Row3 $ws 6 "Elementary" "TestItem_AndSplit:0" "Left1"
Set-Row3 $ws 7 "Elementary" "TestItem_AndSplit:0" "Left2"
Set-Row3 $ws 8 "BlockEnd" "" ""
Set-Row3 $ws 9 "Block" "" ""
Set-Row3 $ws 10 "Elementary" "TestItem_AndSplit:0" "Right1"
Set-Row3 $ws 11 "Elementary" "TestItem_AndSplit:0" "Right2"
Set-Row3 $ws 12 "End" "" ""
Set-Row3 $ws 13 "End" "" ""
Set-Row3 $ws 14 "Elementary" "TestItem_Last:0" "Last"

# ---------------------------------------------------------------------
# Sheet: NestedAndSplit
# ---------------------------------------------------------------------
$ws = New-LayoutSheet "NestedAndSplit" 15

Set-Row3 $ws 1 "layout" $null $null
Set-Row3 $ws 2 "class" "activityReference" "name"
Set-Row3 $ws 3 "AndSplit" "" ""
Set-Row3 $ws 4 "Block" "" ""
Set-Row3 $ws 5 "Elementary" "TestItem_AndSplit:0" "Left"
Set-Row3 $ws 6 "End" "" ""
Set-Row3 $ws 7 "AndSplit" "" ""
Set-Row3 $ws 8 "Block" "" ""
Set-Row3 $ws 9 "Elementary" "TestItem_AndSplit:0" "Right1"
Set-Row3 $ws 10 "End" "" ""
Set-Row3 $ws 11 "Block" "" ""
Set-Row3 $ws 12 "Elementary" "TestItem_AndSplit:0" "Right2"
Set-Row3 $ws 13 "End" "" ""
Set-Row3 $ws 14 "End" "" ""
Set-Row3 $ws 15 "End" "" ""

# ---------------------------------------------------------------------
# Sheet: AndSplitWithLoop
# ---------------------------------------------------------------------
$ws = New-LayoutSheet "AndSplitWithLoop" 10

Set-Row3 $ws 1 "layout" $null $null
Set-Row3 $ws 2 "class" "activityReference" "name"
Set-Row3 $ws 3 "AndSplit" "" ""
Set-Row3 $ws 4 "Loop" "" ""
Set-Row3 $ws 5 "Elementary" "TestItem_AndSplit:0" "Loop"
Set-Row3 $ws 6 "LoopEnd" "" ""
Set-Row3 $ws 7 "Block" "" ""
Set-Row3 $ws 8 "Elementary" "TestItem_AndSplit:0" "Right"
Set-Row3 $ws 9 "End" "" ""
Set-Row3 $ws 10 "End" "" ""

# Make the last added sheet ("AndSplitWithLoop") the active / selected tab,
# matching the final state recorded for this workbook.
$ws.Activate()
